# Adds a new "2022" data column (column S) to the 3.b.1 vaccine-coverage
# sheet, mirroring the formatting already used for the neighbouring
# "2021" column (column R), and fixes up the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Most S-column cells should simply inherit the formatting that the
#    same row already uses in column R (the previous "year" column).
#    Copy R -> S formatting for every row that needs a new S cell.
# ---------------------------------------------------------------------
$rowsMirroringR = @(3,4,5,6,7,8,9,10,11,12,13,14,15,16,19,20,21,22,23,24,25,27,28,29,30,31,32,33,34,35,36)
foreach ($r in $rowsMirroringR) {
    $ws.Range("R$r").Copy() | Out-Null
    $ws.Range("S$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# A few rows (17, 18, 26, 37) pick up the formatting used one row above
# column Q/R's own "first data row" styling instead of their own row's
# R-column style - mirror that explicitly.
$ws.Range("Q6").Copy() | Out-Null
$ws.Range("S17").PasteSpecial(-4122) | Out-Null

$ws.Range("Q7").Copy() | Out-Null
$ws.Range("S18").PasteSpecial(-4122) | Out-Null

$ws.Range("Q15").Copy() | Out-Null
$ws.Range("S26").PasteSpecial(-4122) | Out-Null

$ws.Range("Q7").Copy() | Out-Null
$ws.Range("S37").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Write the 2022 values into column S.
# ---------------------------------------------------------------------
$ws.Range("S4").Value = 2022

$dataRows = @(6,7,8,9,10,11,12,13,14,15,17,18,19,20,21,22,23,24,25,26,28,29,30,31,32,33,34,35,36,37)
$dataVals = @(
    91.320113549242663,
    95.532963647566234,
    91.979142449101602,
    97.11736444749485,
    95.22197889707347,
    95.83359340865114,
    91.694814226107695,
    92.720266061341917,
    78.590540307267389,
    88.700629650829995,
    96.389078828315476,
    95.902649414664197,
    95.851738682785879,
    102.7568062228323,
    99.681465259804895,
    99.066849759690413,
    99.105901053049877,
    100.16892783614,
    87.061971344726402,
    95.855752718946761,
    92.843773094907561,
    95.656192236598898,
    91.580590521106643,
    97.643559826126747,
    95.022894684451515,
    96.0270894451033,
    90.224518180011927,
    94.65237818719315,
    86.721042637666145,
    89.112764739553512
)

for ($i = 0; $i -lt $dataRows.Length; $i++) {
    $ws.Range("S$($dataRows[$i])").Value = $dataVals[$i]
}

# ---------------------------------------------------------------------
# 3) Row 38 (the thick bottom rule under the table) gains a new,
#    previously-unused style: same font/fill as the rest of the row,
#    plus the medium black bottom border that the thick-bottom rule
#    uses elsewhere (borderId 1), and no special alignment.
# ---------------------------------------------------------------------
$s38 = $ws.Range("S38")
$s38.HorizontalAlignment = 1       # xlGeneral
$s38.VerticalAlignment = -4107     # xlBottom (Excel default - no <alignment> emitted)
$s38.Borders.Item(9).LineStyle = 1
$s38.Borders.Item(9).Weight = -4138
$s38.Borders.Item(9).ColorIndex = 1

# ---------------------------------------------------------------------
# 4) Selection moves to the new header cell.
# ---------------------------------------------------------------------
$ws.Range("T4").Select() | Out-Null
